$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 298.42856
$ws.Range("J41").Value = 272
$ws.Range("L41").Value = 272
$ws.Range("N41").Value = -1152
$ws.Range("H70").Value = 3821.5186
$ws.Range("J70").Value = 4124.2085
$ws.Range("L70").Value = 12372.6255
$ws.Range("N70").Value = -12912.6255
$ws.Range("H73").Value = 3821.5186
$ws.Range("J73").Value = 4124.2085
$ws.Range("L73").Value = 12372.6255
$ws.Range("N73").Value = -14244.6255
$ws.Range("H76").Value = 14000
$ws.Range("I76").Value = 12000
$ws.Range("K76").Value = 12000
$ws.Range("M76").Value = -11685
$ws.Range("H79").Value = 14000
$ws.Range("I79").Value = 12000
$ws.Range("K79").Value = 12000
$ws.Range("M79").Value = -10908
$ws.Range("H80").Value = 441.69232
$ws.Range("J80").Value = 382.83334
$ws.Range("L80").Value = 1148.50002
$ws.Range("N80").Value = -3144.50002
$ws.Range("H83").Value = 441.69232
$ws.Range("J83").Value = 382.83334
$ws.Range("L83").Value = 3445.50006
$ws.Range("N83").Value = -13429.50006
$ws.Range("H96").Value = 496.54544
$ws.Range("I96").Value = 536.3
$ws.Range("J96").Value = 99
$ws.Range("K96").Value = 1608.9
$ws.Range("L96").Value = 297
$ws.Range("M96").Value = -235.8999999999999
$ws.Range("N96").Value = -3043
$ws.Range("H112").Value = 3180.077
$ws.Range("J112").Value = 3397.2727
$ws.Range("L112").Value = 10191.8181
$ws.Range("N112").Value = -12407.8181
$ws.Range("H132").Value = 35766.332
$ws.Range("I132").Value = 35766.332
$ws.Range("K132").Value = 107298.996
$ws.Range("M132").Value = -104768.996
$ws.Range("H138").Value = 50002450
$ws.Range("J138").Value = 111113050
$ws.Range("L138").Value = 333339150
$ws.Range("N138").Value = -333349430

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 8000
$ws.Range("J11").Value = 8000
$ws.Range("L11").Value = 8000
$ws.Range("N11").Value = -8288
$ws.Range("H23").Value = 49666.668
$ws.Range("J23").Value = 49500
$ws.Range("L23").Value = 49500
$ws.Range("N23").Value = -50018
$ws.Range("H35").Value = 14500
$ws.Range("J35").Value = 27000
$ws.Range("L35").Value = 27000
$ws.Range("N35").Value = -27812
$ws.Range("H45").Value = 7641.6665
$ws.Range("I45").Value = 4591.6665
$ws.Range("K45").Value = 4591.6665
$ws.Range("M45").Value = -4214.6665
$ws.Range("H102").Value = 4624.263
$ws.Range("I102").Value = 3756.5293
$ws.Range("J102").Value = 12000
$ws.Range("K102").Value = 3756.5293
$ws.Range("L102").Value = 12000
$ws.Range("M102").Value = -2134.5293
$ws.Range("N102").Value = -15244
$ws.Range("H132").Value = 3646.4883
$ws.Range("I132").Value = 2732.7778
$ws.Range("J132").Value = 8345.571
$ws.Range("K132").Value = 8198.3334
$ws.Range("L132").Value = 25036.713
$ws.Range("M132").Value = -5668.3334
$ws.Range("N132").Value = -30096.713

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5317.5625
$ws.Range("J20").Value = 6213.7144
$ws.Range("L20").Value = 6213.7144
$ws.Range("N20").Value = -6707.7144
$ws.Range("H105").Value = 4100
$ws.Range("I105").Value = 2850
$ws.Range("J105").Value = 5975
$ws.Range("K105").Value = 2850
$ws.Range("L105").Value = 5975
$ws.Range("M105").Value = -1103
$ws.Range("N105").Value = -9469

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 2937.5
$ws.Range("I12").Value = 3250
$ws.Range("J12").Value = 2000
$ws.Range("K12").Value = 3250
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = -3080
$ws.Range("N12").Value = -2340
$ws.Range("H22").Value = 762.25
$ws.Range("J22").Value = 899.6667
$ws.Range("L22").Value = 899.6667
$ws.Range("N22").Value = -1599.6667
$ws.Range("H58").Value = 9104.857
$ws.Range("I58").Value = 5979.8335
$ws.Range("K58").Value = 5979.8335
$ws.Range("M58").Value = -5776.8335
$ws.Range("H132").Value = 2777.5334
$ws.Range("I132").Value = 2833.0715
$ws.Range("K132").Value = 8499.2145
$ws.Range("M132").Value = -5969.2145
$ws.Range("H133").Value = 55882.332
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 55882.332
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 55882.332
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -60942.332
$ws.Range("H136").Value = 9104.857
$ws.Range("I136").Value = 5979.8335
$ws.Range("K136").Value = 17939.5005
$ws.Range("M136").Value = -15389.5005

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 2583808
$ws.Range("J9").Value = 528
$ws.Range("L9").Value = 1584
$ws.Range("N9").Value = -2032
$ws.Range("H55").Value = 700
$ws.Range("I55").Value = 700
$ws.Range("K55").Value = 2100
$ws.Range("M55").Value = -1923
$ws.Range("H56").Value = 5000
$ws.Range("I56").Value = 5000
$ws.Range("K56").Value = 5000
$ws.Range("M56").Value = -4470
$ws.Range("H68").Value = 3640.5
$ws.Range("I68").Value = 742.7143
$ws.Range("K68").Value = 2228.1429
$ws.Range("M68").Value = -1417.1429
$ws.Range("H71").Value = 3640.5
$ws.Range("I71").Value = 742.7143
$ws.Range("K71").Value = 6684.428699999999
$ws.Range("M71").Value = -2628.428699999999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 198.75
$ws.Range("I17").Value = 198.75
$ws.Range("K17").Value = 198.75
$ws.Range("M17").Value = -30.75
$ws.Range("H102").Value = 2436.4443
$ws.Range("I102").Value = 1834.75
$ws.Range("K102").Value = 1834.75
$ws.Range("M102").Value = -212.75
$ws.Range("H132").Value = 3256.5862
$ws.Range("I132").Value = 2807.261
$ws.Range("K132").Value = 8421.782999999999
$ws.Range("M132").Value = -5891.782999999999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 5000
$ws.Range("I13").Value = 1000
$ws.Range("J13").Value = 9000
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 9000
$ws.Range("M13").Value = -860
$ws.Range("N13").Value = -9280
$ws.Range("H55").Value = 639.06665
$ws.Range("I55").Value = 880.7143
$ws.Range("K55").Value = 880.7143
$ws.Range("M55").Value = -707.7143
$ws.Range("H82").Value = 39999.5
$ws.Range("I82").Value = 39999.5
$ws.Range("K82").Value = 39999.5
$ws.Range("M82").Value = -39638.5
$ws.Range("H85").Value = 39999.5
$ws.Range("I85").Value = 39999.5
$ws.Range("K85").Value = 39999.5
$ws.Range("M85").Value = -38751.5
$ws.Range("H130").Value = 36994.8
$ws.Range("J130").Value = 36994.8
$ws.Range("L130").Value = 36994.8
$ws.Range("N130").Value = -47034.8
$ws.Range("H136").Value = 5444.4
$ws.Range("I136").Value = 4853.3335
$ws.Range("K136").Value = 14560.0005
$ws.Range("M136").Value = -12010.0005

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 33333
$ws.Range("J68").Value = 33333
$ws.Range("L68").Value = 33333
$ws.Range("N68").Value = -34955
$ws.Range("H71").Value = 33333
$ws.Range("J71").Value = 33333
$ws.Range("L71").Value = 99999
$ws.Range("N71").Value = -108111
$ws.Range("H81").Value = 7266
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 7266
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H92").Value = 45000
$ws.Range("J92").Value = 45000
$ws.Range("L92").Value = 45000
$ws.Range("N92").Value = -49992
$ws.Range("H100").Value = 1278.5385
$ws.Range("I100").Value = 1282.7
$ws.Range("K100").Value = 2565.4
$ws.Range("M100").Value = -2024.4
$ws.Range("H132").Value = 3343.8076
$ws.Range("I132").Value = 3501.9048
$ws.Range("K132").Value = 10505.7144
$ws.Range("M132").Value = -7975.714399999999
